$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.036415338516235
$ws.Range("B1").Value = 2.271584510803223
$ws.Range("C1").Value = 4.465034484863281
$ws.Range("D1").Value = 1.10789954662323
$ws.Range("E1").Value = 1.251142144203186
